$d = $word.ActiveDocument

# ---- Chunk A: replace paragraphs 1-13 (title through Kronecker/exponentials block) ----
$pStart = $d.Paragraphs(1)
$pEnd = $d.Paragraphs(13)
$rangeA = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$chunkA = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>Benchmarks</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> (solid effect</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>, laptop</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Calculating 1E6 matrix products: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Matlab (s): </w:t></w:r><w:r><w:t>1.34</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_Hlk518317474"/><w:r><w:t>Python + NumPy</w:t></w:r><w:r><w:t>, Enthought</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">(s): </w:t></w:r><w:r><w:t>2.23</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Python + Numpy</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>Intel (s): 2.</w:t></w:r><w:r><w:t>27</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Python + Numpy, </w:t></w:r><w:r><w:t>Anaconda</w:t></w:r><w:r><w:t xml:space="preserve"> (s): 2.</w:t></w:r><w:r><w:t>38</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Python + F2PY (s):</w:t></w:r><w:r><w:t xml:space="preserve"> 1.90</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Fortran</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gfortran</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (s)</w:t></w:r><w:r><w:t>: 1.88</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Fortran, </w:t></w:r><w:r><w:t>Intel</w:t></w:r><w:r><w:t xml:space="preserve"> (s): </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:t>Try Fortran with Intel MKL then on Ubuntu machine. If still slower than Matlab</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> post on Stackoverflow for help. This does however explain why Matlab program is much faster than Python.</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> It is curious that Fortran is slower than Matlab, presumably due to high optimisation of Matlab matrix multiplication?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/></w:rPr></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:t>This would mean there is no point optimising Python code, performance will never exceed Matlab?</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t xml:space="preserve">Calculating </w:t></w:r><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t>1E5</w:t></w:r><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t xml:space="preserve"> Kronecker products: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Matlab (s): </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>4.18</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Python + NumPy (s): 7.33, 7.34, 7.41</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t xml:space="preserve">Calculating </w:t></w:r><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t>1E4</w:t></w:r><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t xml:space="preserve"> matrix exponentials: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Matlab (s): </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>2.02</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Python + NumPy (s): 6.77, 6.70, 6.73</w:t></w:r></w:p>
'@

$rangeA.InsertXML($chunkA)

# ---- Chunk B: merge runs in the "Python + F2PY (s): 0.070, 0.070, 0.070" paragraph ----
# (Hamiltonian section, now at paragraph 24 after Chunk A inserted 6 extra paragraphs)
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "Python + F2PY (s):*0.070*") {
        $targetPara = $d.Paragraphs($i)
        break
    }
}

$chunkB = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Python + </w:t></w:r><w:r><w:t>F2PY</w:t></w:r><w:r><w:t xml:space="preserve"> (s): </w:t></w:r><w:r><w:t>0.070, 0.070, 0.070</w:t></w:r></w:p>
'@

$targetPara.Range.InsertXML($chunkB)

Write-Output "Done"
